$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 80 ("Added semester number to display (JG - 2018/02/03)")
# which pushes every row from the old row 80 onward down by one.
$ws.Rows(80).Insert()
$ws.Range("D80").Value = "Added semester number to display (JG - 2018/02/03)"
$ws.Range("D80").Font.Bold = $false
$ws.Range("D80").HorizontalAlignment = -4131

# Resize columns A, B, C, G (column widths changed by the author)
$ws.Columns("A").ColumnWidth = 19.877604166666668
$ws.Columns("B").ColumnWidth = 17.736979166666668
$ws.Columns("C").ColumnWidth = 9.307291666666666
$ws.Columns("G").ColumnWidth = 19.451822916666668

# Update the view: scroll position reset and new active selection
$ws.Activate() | Out-Null
$ws.Range("L80").Select() | Out-Null
